# Updates cryptos list prices / 1h-volume percentages for the scraped
# coinranking.com table, plus swaps the Mantle <-> NEARProtocol rows back
# into their "correct" rank order (30 = NEARProtocol, 31 = Mantle).
#
# Price cells in column D are stored as literal text (they contain things
# like "0.999", "1.00", "65.609.59", "0.0₃0638" -- thousand-dot formatted
# numbers, values with significant trailing/leading zeros, and even a
# subscript-zero notation) so a plain `.Value = "1.00"` assignment would
# get auto-coerced by Excel into the *number* 1, destroying the intended
# text. Forcing NumberFormat="@" (Text) before the assignment keeps the
# literal string, and resetting `.Style` back to 'Normal' afterwards keeps
# the cell's visual style identical to the untouched cells around it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '65.609.59'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +0.94%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.189.48'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +0.89%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('D4').Style = 'Normal'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '594.15'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +3.99%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '150.58'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.26%  '
$ws.Range('E7').Value = '  -0.13%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '3.187.51'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +1.04%  '
$ws.Range('E9').Value = '  +1.61%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.161'
$ws.Range('D10').Style = 'Normal'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '6.15'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.83%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.509'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +1.00%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000270'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -0.86%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '38.11'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.03%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.712.31'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.90%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '65.825.09'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +1.11%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '7.29'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +1.46%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '3.197.47'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +1.03%  '
$ws.Range('E19').Value = '  +0.33%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '511.33'
$ws.Range('D20').Style = 'Normal'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '15.87'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +6.71%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.730'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.30%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '15.24'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -3.67%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '7.93'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +1.38%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '85.38'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.71%  '
$ws.Range('E26').Value = '  +0.00%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.23'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +1.48%  '
$ws.Range('E28').Value = '  +3.54%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.23'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +1.55%  '
$ws.Range('E30').Value = '  +3.01%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '28.00'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.14%  '
$ws.Range('B32').Value = 'NEARProtocol'
$ws.Range('C32').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '6.66'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +5.35%  '
$ws.Range('B33').Value = 'Mantle'
$ws.Range('C33').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.23'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.54%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.00'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.16%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '6.59'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.56%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '55.38'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.37%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.0914'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +3.98%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '480.96'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +1.10%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0424'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.83%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.98'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -3.40%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '8.95'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +3.88%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.019.87'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -3.10%  '
$ws.Range('E43').Value = '  -2.08%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.289'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.24%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.45'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -2.90%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0₃0638'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +8.98%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '29.00'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.98%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.999'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.03%  '
$ws.Range('E49').Value = '  +0.53%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.29'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.05%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '119.99'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -2.32%  '
